$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Liste")

$ws.Range("F3").Value = "U3-4"
$ws.Range("F5").Value = "U3-110"
$ws.Range("F8").Value = "U3-110"
$ws.Range("F10").Value = "U3-Amphi"
$ws.Range("F13").Value = "U3-Amphi"
$ws.Range("F15").Value = "U3-110"
$ws.Range("F16").Value = "U3-110"
$ws.Range("F17").Value = "U3-109"
$ws.Range("F20").Value = "U3-Amphi"
